$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the existing row 58, pushing
# all subsequent rows (58-170) down by one (to 59-171). Insert a blank row
# at position 58 first so the rest of the sheet shifts down automatically.
$ws.Rows(58).Insert()

# Populate the newly inserted row 58 with the new record's data.
$ws.Cells.Item(58, 1).Value = 3
$ws.Cells.Item(58, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(58, 3).Value = "Coquimbo"
$ws.Cells.Item(58, 4).Value = 44662
$ws.Cells.Item(58, 5).Value = 5
$ws.Cells.Item(58, 6).Value = 100112052
$ws.Cells.Item(58, 7).Value = "Albahaca"
$ws.Cells.Item(58, 8).Value = "Sin especificar"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 80
$ws.Cells.Item(58, 11).Value = 4000
$ws.Cells.Item(58, 12).Value = 4000
$ws.Cells.Item(58, 13).Value = 4000
$ws.Cells.Item(58, 14).Value = "`$/docena de matas"
$ws.Cells.Item(58, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(58, 16).Value = 667
$ws.Cells.Item(58, 17).Value = 6
$ws.Cells.Item(58, 18).Value = "Hortaliza"
